# Estadisticos Segundo Parcial 23 Mayo
# Updates second-partial stats, the final stats (which mirror 2P's
# Reprobados/Aprobados/Por_Apro but keep their own Promedio), and replaces
# the "Rescatables" (students who need a make-up exam) roster.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "Estadisticos 2P" — Blancos/Reprobados/Aprobados/Por_Apro/Promedio
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Cells.Item(2,4).Value = 0
$ws2.Cells.Item(2,5).Value = 7
$ws2.Cells.Item(2,6).Value = 29
$ws2.Cells.Item(2,7).Value = 80.56
$ws2.Cells.Item(2,8).Value = 6.1

$ws2.Cells.Item(3,4).Value = 0
$ws2.Cells.Item(3,5).Value = 5
$ws2.Cells.Item(3,6).Value = 14
$ws2.Cells.Item(3,7).Value = 73.68000000000001
$ws2.Cells.Item(3,8).Value = 6.6

$ws2.Cells.Item(4,4).Value = 0
$ws2.Cells.Item(4,5).Value = 4
$ws2.Cells.Item(4,6).Value = 16
$ws2.Cells.Item(4,7).Value = 80
$ws2.Cells.Item(4,8).Value = 6

# ---------------------------------------------------------------
# "Estadisticos Final" — same Reprobados/Aprobados/Por_Apro as 2P,
# but its own Promedio values.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Cells.Item(2,5).Value = 7
$ws3.Cells.Item(2,6).Value = 29
$ws3.Cells.Item(2,7).Value = 80.56
$ws3.Cells.Item(2,8).Value = 6.8

$ws3.Cells.Item(3,5).Value = 5
$ws3.Cells.Item(3,6).Value = 14
$ws3.Cells.Item(3,7).Value = 73.68000000000001
$ws3.Cells.Item(3,8).Value = 6.7

$ws3.Cells.Item(4,5).Value = 4
$ws3.Cells.Item(4,6).Value = 16
$ws3.Cells.Item(4,7).Value = 80
$ws3.Cells.Item(4,8).Value = 6.5

# ---------------------------------------------------------------
# "Rescatables" — roster of students that still need a make-up
# exam. The whole table is replaced (9 rows -> 14 rows).
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$ws4.Cells.Item(2,1).Value = 24330051920092
$ws4.Cells.Item(2,2).Value = "APARICIO"
$ws4.Cells.Item(2,3).Value = "OFICIAL"
$ws4.Cells.Item(2,4).Value = "VICTOR YAEL"
$ws4.Cells.Item(2,5).Value = "DISEÑA INSTALACIONES ELÉCTRICAS"
$ws4.Cells.Item(2,6).Value = "2AEV"
$ws4.Cells.Item(2,7).Value = 4

$ws4.Cells.Item(3,1).Value = 24330051920122
$ws4.Cells.Item(3,2).Value = "RIOS"
$ws4.Cells.Item(3,3).Value = "ZEPAHUA"
$ws4.Cells.Item(3,4).Value = "UZIEL"
$ws4.Cells.Item(3,5).Value = "DISEÑA INSTALACIONES ELÉCTRICAS"
$ws4.Cells.Item(3,6).Value = "2AEV"
$ws4.Cells.Item(3,7).Value = 4

$ws4.Cells.Item(4,1).Value = 24330051920315
$ws4.Cells.Item(4,2).Value = "VENTURA"
$ws4.Cells.Item(4,3).Value = "ZEPEDA"
$ws4.Cells.Item(4,4).Value = "CARLOS ARGEL"
$ws4.Cells.Item(4,5).Value = "DISEÑA INSTALACIONES ELÉCTRICAS"
$ws4.Cells.Item(4,6).Value = "2AEV"
$ws4.Cells.Item(4,7).Value = 4

$ws4.Cells.Item(5,1).Value = 23330051920211
$ws4.Cells.Item(5,2).Value = "VAZQUEZ"
$ws4.Cells.Item(5,3).Value = "CARRILLO"
$ws4.Cells.Item(5,4).Value = "DIEGO ARMANDO"
$ws4.Cells.Item(5,5).Value = "PROGRAMA Y CONECTA CONTROLADORES LÓGICOS PROGRAMABLES (PLC´S)"
$ws4.Cells.Item(5,6).Value = "4AEV"
$ws4.Cells.Item(5,7).Value = 4

$ws4.Cells.Item(6,1).Value = 21330051920007
$ws4.Cells.Item(6,2).Value = "COBOS"
$ws4.Cells.Item(6,3).Value = "NOLASCO"
$ws4.Cells.Item(6,4).Value = "YOLET"
$ws4.Cells.Item(6,5).Value = "REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA"
$ws4.Cells.Item(6,6).Value = "6AEV"
$ws4.Cells.Item(6,7).Value = 4

$ws4.Cells.Item(7,1).Value = 22330051920413
$ws4.Cells.Item(7,2).Value = "LOBATO"
$ws4.Cells.Item(7,3).Value = "ANTONIO"
$ws4.Cells.Item(7,4).Value = "FABIAN ALEJANDRO"
$ws4.Cells.Item(7,5).Value = "REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA"
$ws4.Cells.Item(7,6).Value = "6AEV"
$ws4.Cells.Item(7,7).Value = 4

$ws4.Cells.Item(8,1).Value = 24330051920090
$ws4.Cells.Item(8,2).Value = "ANTONIO"
$ws4.Cells.Item(8,3).Value = "LOPEZ"
$ws4.Cells.Item(8,4).Value = "SERGIO GISELL"
$ws4.Cells.Item(8,5).Value = "DISEÑA INSTALACIONES ELÉCTRICAS"
$ws4.Cells.Item(8,6).Value = "2AEV"
$ws4.Cells.Item(8,7).Value = 3

$ws4.Cells.Item(9,1).Value = 23330051920036
$ws4.Cells.Item(9,2).Value = "HERNANDEZ"
$ws4.Cells.Item(9,3).Value = "DOLORES"
$ws4.Cells.Item(9,4).Value = "GONZALO"
$ws4.Cells.Item(9,5).Value = "DISEÑA INSTALACIONES ELÉCTRICAS"
$ws4.Cells.Item(9,6).Value = "2AEV"
$ws4.Cells.Item(9,7).Value = 3

$ws4.Cells.Item(10,1).Value = 24330051920107
$ws4.Cells.Item(10,2).Value = "PELLICO"
$ws4.Cells.Item(10,3).Value = "SANCHEZ"
$ws4.Cells.Item(10,4).Value = "MIRANDA ALIZEET"
$ws4.Cells.Item(10,5).Value = "DISEÑA INSTALACIONES ELÉCTRICAS"
$ws4.Cells.Item(10,6).Value = "2AEV"
$ws4.Cells.Item(10,7).Value = 3

$ws4.Cells.Item(11,1).Value = 22330051920389
$ws4.Cells.Item(11,2).Value = "FLORES"
$ws4.Cells.Item(11,3).Value = "LAGUNA"
$ws4.Cells.Item(11,4).Value = "JOSE ANTONIO"
$ws4.Cells.Item(11,5).Value = "PROGRAMA Y CONECTA CONTROLADORES LÓGICOS PROGRAMABLES (PLC´S)"
$ws4.Cells.Item(11,6).Value = "4AEV"
$ws4.Cells.Item(11,7).Value = 3

$ws4.Cells.Item(12,1).Value = 23330051920332
$ws4.Cells.Item(12,2).Value = "RODRIGUEZ"
$ws4.Cells.Item(12,3).Value = "SUAREZ"
$ws4.Cells.Item(12,4).Value = "SERGIO JOSUE"
$ws4.Cells.Item(12,5).Value = "PROGRAMA Y CONECTA CONTROLADORES LÓGICOS PROGRAMABLES (PLC´S)"
$ws4.Cells.Item(12,6).Value = "4AEV"
$ws4.Cells.Item(12,7).Value = 3

$ws4.Cells.Item(13,1).Value = 22330051920177
$ws4.Cells.Item(13,2).Value = "CAMPOS"
$ws4.Cells.Item(13,3).Value = "CABRERA"
$ws4.Cells.Item(13,4).Value = "MARCO"
$ws4.Cells.Item(13,5).Value = "REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA"
$ws4.Cells.Item(13,6).Value = "6AEV"
$ws4.Cells.Item(13,7).Value = 3

$ws4.Cells.Item(14,1).Value = 22330051920193
$ws4.Cells.Item(14,2).Value = "MOLINA"
$ws4.Cells.Item(14,3).Value = "DE JESUS"
$ws4.Cells.Item(14,4).Value = "VICTOR MANUEL"
$ws4.Cells.Item(14,5).Value = "REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA"
$ws4.Cells.Item(14,6).Value = "6AEV"
$ws4.Cells.Item(14,7).Value = 2
